$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the dataset. Insert a fresh row at
# position 11 (shifting the existing rows 11-18 down to 12-19) and fill
# it with the new observation.
$ws.Rows.Item(11).Insert()

$ws.Cells.Item(11, 1).Value = 7
$ws.Cells.Item(11, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(11, 3).Value = "Ñuble"
$ws.Cells.Item(11, 4).Value = 45176
$ws.Cells.Item(11, 5).Value = 16
$ws.Cells.Item(11, 6).Value = "Fruta"
$ws.Cells.Item(11, 7).Value = 100107
$ws.Cells.Item(11, 8).Value = "Otros"
$ws.Cells.Item(11, 9).Value = 100107002
$ws.Cells.Item(11, 10).Value = "Chirimoya"
$ws.Cells.Item(11, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 30
$ws.Cells.Item(11, 14).Value = 22000
$ws.Cells.Item(11, 15).Value = 22000
$ws.Cells.Item(11, 16).Value = 22000
$ws.Cells.Item(11, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(11, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(11, 19).Value = 2200
$ws.Cells.Item(11, 20).Value = 10
